$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header date text (kept as plain text, matching original cell formatting)
$ws.Range("B1").Value = "25/03/2023"

# Update hourly values
$ws.Range("B2").Value = 78
$ws.Range("B3").Value = 167
$ws.Range("B4").Value = 138
$ws.Range("B5").Value = 125
$ws.Range("B6").Value = 116
$ws.Range("B7").Value = 99
$ws.Range("B8").Value = 98
$ws.Range("B9").Value = 127
$ws.Range("B10").Value = 105
$ws.Range("B11").Value = 105
$ws.Range("B12").Value = 91
$ws.Range("B13").Value = 82
$ws.Range("B14").Value = 68
$ws.Range("B15").Value = 27
$ws.Range("B16").Value = 27
$ws.Range("B17").Value = 22
